# Apply the hexsticker.pptx edit:
#  - Slide 1: resize/reposition "Picture 4" (the faded/transparent hex outline
#    image) and delete "Picture 7" (the small standalone logo picture).
#  - Slide 2: shift "Picture 7" (the big sticker screenshot) to the right and
#    delete the "Group 5" group (background rectangle + cropped logo photo).

$p = $ppt.ActivePresentation

# ---- Slide 1 -------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# "Picture 4": off 4213904,1299261 ext 2083713x1733650 -> off 3924193,1299261
# ext 2687876x2236314  (EMU / 12700 = points)
$pic4 = $s1.Shapes.Item("Picture 4")
$pic4.Left   = 3924193 / 12700
$pic4.Top    = 1299261 / 12700
$pic4.Width  = 2687876 / 12700
$pic4.Height = 2236314 / 12700

# "Picture 7" (small logo pic, off 4899056,2838073 ext 723146x723146) removed.
$s1.Shapes.Item("Picture 7").Delete()

# ---- Slide 2 -------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# "Picture 7": off 170684,545823 -> off 3466145,545823 (ext unchanged)
$pic7b = $s2.Shapes.Item("Picture 7")
$pic7b.Left = 3466145 / 12700

# "Group 5" (Rectangle 4 + Rectangle 3) removed entirely.
$s2.Shapes.Item("Group 5").Delete()
